# ETI ToDoList #37 - Update test case documentation (row 5 retest results)
# and refresh the viewport/selection to match the reviewed state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 ("Test displaying of individual to-do history for each user") was
# retested: it now passed, with a justification note instead of the
# previous "Fail" marker.
$ws.Range("F5").Value = "Same as expected outcome."
$ws.Range("G5").Value = "Pass"

# Move the on-screen selection/viewport to the cell that was being
# reviewed (G4:G5, scrolled so column D is the first visible column).
$ws.Range("G4:G5").Select()
